# Apply the "ome" prefix rename across the workbook.
# Previously the base (blank-name) RDF prefix was written as ":" — this
# commit gives it an explicit name "ome" in the @prefix sheet, and every
# cell that previously held a bare ":xxx" / ":Xxx" CURIE is rewritten to
# use the new "ome:" prefix instead.

$wb = $excel.ActiveWorkbook

# 1) @prefix sheet: name the base prefix "ome" (was blank).
$wsPrefix = $wb.Worksheets.Item("@prefix")
$wsPrefix.Range("A1").Value = "ome"

# 2) Image sheet: CURIEs that used the bare ":" prefix now use "ome:".
$wsImage = $wb.Worksheets.Item("Image")
$wsImage.Range("E3").Value = "ome:pixels"
$wsImage.Range("F3").Value = "ome:acquisitionDate"
$wsImage.Range("B4").Value = "ome:Image"
$wsImage.Range("E4").Value = "ome:Pixels"
$wsImage.Range("F4").Select()

# 3) Pixels sheet.
$wsPixels = $wb.Worksheets.Item("Pixels")
$wsPixels.Range("D3").Value = "ome:pixelType"
$wsPixels.Range("E3").Value = "ome:dimensionOrder"
$wsPixels.Range("F3").Value = "ome:sizeC"
$wsPixels.Range("G3").Value = "ome:sizeT"
$wsPixels.Range("H3").Value = "ome:sizeX"
$wsPixels.Range("I3").Value = "ome:sizeY"
$wsPixels.Range("J3").Value = "ome:sizeZ"
$wsPixels.Range("K3").Value = "ome:channel"
$wsPixels.Range("L3").Value = "ome:binData"
$wsPixels.Range("B4").Value = "ome:Pixels"
$wsPixels.Range("D4").Value = "ome:PixelType"
$wsPixels.Range("E4").Value = "ome:DimensionOrder"
$wsPixels.Range("K4").Value = "ome:Channel"
$wsPixels.Range("L4").Value = "ome:BinData"

# 4) Channel sheet.
$wsChannel = $wb.Worksheets.Item("Channel")
$wsChannel.Range("D3").Value = "ome:color"
$wsChannel.Range("B4").Value = "ome:Channel"
$wsChannel.Range("D4").Value = "ome:Color"

# 5) Color sheet.
$wsColor = $wb.Worksheets.Item("Color")
$wsColor.Range("C3").Value = "ome:length"
$wsColor.Range("B4").Value = "ome:Color"

# 6) Binary_Data sheet.
$wsBin = $wb.Worksheets.Item("Binary_Data")
$wsBin.Range("C3").Value = "ome:bigEndian"
$wsBin.Range("D3").Value = "ome:data"
$wsBin.Range("E3").Value = "ome:length"
$wsBin.Range("B4").Value = "ome:BinData"

# Binary_Data ends up the active tab after the edit session.
$wsBin.Activate()
